$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.285.62"
$ws.Range("E2").Value = "'  +2.10%  "
$ws.Range("D3").Value = "'2.098.31"
$ws.Range("E3").Value = "'  +0.12%  "
$ws.Range("E4").Value = "'  -0.56%  "
$ws.Range("D5").Value = "'342.69"
$ws.Range("E5").Value = "'  -0.22%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  -0.54%  "
$ws.Range("D7").Value = "'0.5267"
$ws.Range("E7").Value = "'  +2.29%  "
$ws.Range("D8").Value = "'0.4405"
$ws.Range("E8").Value = "'  +0.74%  "
$ws.Range("D9").Value = "'55.03"
$ws.Range("E9").Value = "'  +3.00%  "
$ws.Range("D10").Value = "'0.09371"
$ws.Range("E10").Value = "'  +2.16%  "
$ws.Range("D11").Value = "'1.177"
$ws.Range("E11").Value = "'  +1.01%  "
$ws.Range("D12").Value = "'24.81"
$ws.Range("E12").Value = "'  +1.21%  "
$ws.Range("D13").Value = "'8.544"
$ws.Range("E13").Value = "'  +4.34%  "
$ws.Range("D14").Value = "'6.891"
$ws.Range("E14").Value = "'  +2.06%  "
$ws.Range("D15").Value = "'2.037.34"
$ws.Range("E15").Value = "'  -2.54%  "
$ws.Range("D16").Value = "'101.45"
$ws.Range("E16").Value = "'  -0.10%  "
$ws.Range("D17").Value = "'0.00001157"
$ws.Range("E17").Value = "'  +0.54%  "
$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = "'  -0.59%  "
$ws.Range("D19").Value = "'21.17"
$ws.Range("E19").Value = "'  +1.20%  "
$ws.Range("D20").Value = "'0.06721"
$ws.Range("E20").Value = "'  +0.72%  "
$ws.Range("D21").Value = "'6.464"
$ws.Range("E21").Value = "'  +4.13%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "'  -0.60%  "
$ws.Range("D23").Value = "'30.304.53"
$ws.Range("E23").Value = "'  +1.96%  "
$ws.Range("D24").Value = "'12.45"
$ws.Range("E24").Value = "'  -0.03%  "
$ws.Range("D25").Value = "'2.321"
$ws.Range("E25").Value = "'  +0.78%  "
$ws.Range("B26").Value = "'EthereumClassic"
$ws.Range("C26").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'21.85"
$ws.Range("E26").Value = "'  -0.16%  "
$ws.Range("B27").Value = "'InternetComputer(DFINITY)"
$ws.Range("C27").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'6.952"
$ws.Range("E27").Value = "'  +10.33%  "
$ws.Range("B28").Value = "'LidoDAOToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.532"
$ws.Range("E28").Value = "'  +1.56%  "
$ws.Range("B29").Value = "'Monero"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'162.22"
$ws.Range("E29").Value = "'  +0.45%  "
$ws.Range("D30").Value = "'133.62"
$ws.Range("E30").Value = "'  +0.09%  "
$ws.Range("D31").Value = "'1.143"
$ws.Range("E31").Value = "'  +1.58%  "
$ws.Range("D32").Value = "'1.672"
$ws.Range("E32").Value = "'  +0.24%  "
$ws.Range("E33").Value = "'  +0.55%  "
$ws.Range("E34").Value = "'  +1.17%  "
$ws.Range("D35").Value = "'3.877"
$ws.Range("E35").Value = "'  -1.59%  "
$ws.Range("D36").Value = "'10.12"
$ws.Range("E36").Value = "'  -2.11%  "
$ws.Range("D37").Value = "'0.02643"
$ws.Range("E37").Value = "'  +2.78%  "
$ws.Range("D38").Value = "'0.06784"
$ws.Range("E38").Value = "'  +1.74%  "
$ws.Range("D39").Value = "'12.70"
$ws.Range("E39").Value = "'  +1.93%  "
$ws.Range("D40").Value = "'0.6968"
$ws.Range("E40").Value = "'  -0.17%  "
$ws.Range("D41").Value = "'1.345"
$ws.Range("E41").Value = "'  +1.15%  "
$ws.Range("D42").Value = "'0.2220"
$ws.Range("E42").Value = "'  +0.13%  "
$ws.Range("D43").Value = "'0.6783"
$ws.Range("E43").Value = "'  -0.46%  "
$ws.Range("D44").Value = "'14.48"
$ws.Range("E44").Value = "'  +1.60%  "
$ws.Range("D45").Value = "'2.332"
$ws.Range("E45").Value = "'  +1.28%  "
$ws.Range("E46").Value = "'  -0.55%  "
$ws.Range("D47").Value = "'1.314"
$ws.Range("E47").Value = "'  +8.29%  "
$ws.Range("D48").Value = "'3.635"
$ws.Range("E48").Value = "'  +0.64%  "
$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000351"
$ws.Range("E49").Value = "'  -0.07%  "
$ws.Range("B50").Value = "'ThetaToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'1.212"
$ws.Range("E50").Value = "'  +6.99%  "
$ws.Range("D51").Value = "'1.215"
$ws.Range("E51").Value = "'  +0.03%  "
